{"js": "// Map of original paragraph text -> replacement text. Using full-paragraph\n// insertText(..., replace) merges each paragraph into a single run, which\n// (as a side effect, matching the real edit) also drops the stray\n// <w:proofErr/> spell-check markers that wrapped \"Matlab\"/\"wavCapture\"/\n// \"shm\"/\"jit\" in the original file.\nconst replacements = {\n  \"To have both of these programs to work you need to have Matlab installed. \":\n    \"To have both programs work, you need to have MATLAB 2013b installed. \",\n  \"Go into the moving window program folder. Select wavCapture code file \":\n    \"Go into the moving window program folder. Select wavCapture code file \",\n  \"Run the program on Matlab and it will open wavCapture.\":\n    \"Run the program on MATLAB and it will open wavCapture.\",\n  \"Evaluate for CPP, D2, and %shm, %jit, and SNR with TF32. \":\n    \"Evaluate for CPP, D2, and %shm, %jit, and SNR with TF32. \",\n  \"Use Excel to rank these parameters (SNR is different than %shm and %jit). Add up the rankings. The one with the lowest value is the one with the least perturbation.\":\n    \"Use Excel to rank these parameters (SNR is different than %shm and %jit). Add up the rankings. The one with the lowest value is the one with the least perturbation.\",\n  \"Open exeFile6 in Matlab.\":\n    \"Open exeFile6 in MATLAB.\",\n  \"Run voiceEvaluate code in Matlab.\":\n    \"Run voiceEvaluate code in MATLAB.\",\n  \"Dr. Lin wrote these programs on Matlab and would be a better person to address questions to. \":\n    \"Dr. Lin wrote these programs on MATLAB and would be a better person to address questions to. \"\n};\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  const replacement = replacements[paragraph.text];\n  if (replacement !== undefined) {\n    paragraph.insertText(replacement, Word.InsertLocation.replace);\n  }\n}\n\n// Remove the leftover \"_GoBack\" bookmark left after \"Moving window program: \".\ncontext.document.deleteBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-DocText($doc, $findText, $replaceText) {\n    $find = $doc.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $replaceText\n    $find.Execute([ref]$findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\n# Rewrite the intro sentence (content change + drops the \"Matlab\" spell-check wrap)\nReplace-DocText $d \"To have both of these programs to work you need to have Matlab installed. \" \"To have both programs work, you need to have MATLAB 2013b installed. \"\n\n# Remove the leftover \"_GoBack\" bookmark after \"Moving window program: \"\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# Drop spell-check wrapping around \"wavCapture\" (text unchanged)\nReplace-DocText $d \"Go into the moving window program folder. Select wavCapture code file \" \"Go into the moving window program folder. Select wavCapture code file \"\n\n# Matlab -> MATLAB, and drop spell-check wrapping around \"wavCapture\"\nReplace-DocText $d \"Run the program on Matlab and it will open wavCapture.\" \"Run the program on MATLAB and it will open wavCapture.\"\n\n# Drop spell-check wrapping around \"shm\" / \"jit\" (text unchanged)\nReplace-DocText $d \"Evaluate for CPP, D2, and %shm, %jit, and SNR with TF32. \" \"Evaluate for CPP, D2, and %shm, %jit, and SNR with TF32. \"\nReplace-DocText $d \"(SNR is different than %shm and %jit). Add up the rankings. The one with the lowest value is the one with the least perturbation.\" \"(SNR is different than %shm and %jit). Add up the rankings. The one with the lowest value is the one with the least perturbation.\"\n\n# Matlab -> MATLAB\nReplace-DocText $d \"Open exeFile6 in Matlab.\" \"Open exeFile6 in MATLAB.\"\nReplace-DocText $d \"Run voiceEvaluate code in Matlab.\" \"Run voiceEvaluate code in MATLAB.\"\nReplace-DocText $d \"Dr. Lin wrote these programs on Matlab and would be a better person to address questions to. \" \"Dr. Lin wrote these programs on MATLAB and would be a better person to address questions to. \"\n"}
